# Apply scheduled market-price refresh to profit-calc sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 8
$ws.Range("H8").Value = 1904.8
$ws.Range("J8").Value = 1.5
$ws.Range("L8").Value = 4.5
$ws.Range("N8").Value = -282.5

# row 74
$ws.Range("H74").Value = 5112.7
$ws.Range("I74").Value = 5112.7
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5112.7
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4176.7
$ws.Range("N74").ClearContents()

# row 77
$ws.Range("H77").Value = 5112.7
$ws.Range("I77").Value = 5112.7
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 25563.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -20883.5
$ws.Range("N77").ClearContents()

# row 96
$ws.Range("H96").Value = 806813.2
$ws.Range("I96").Value = 926.9091
$ws.Range("J96").Value = 2073205.9
$ws.Range("K96").Value = 2780.7273
$ws.Range("L96").Value = 6219617.699999999
$ws.Range("M96").Value = -1407.7273
$ws.Range("N96").Value = -6222363.699999999

# row 97
$ws.Range("H97").Value = 3081.1667
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3081.1667
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 9243.500100000001
$ws.Range("N97").Value = -10235.5001
$ws.Range("M97").ClearContents()

# row 101
$ws.Range("H101").Value = 643.25
$ws.Range("I101").Value = 565.5
$ws.Range("J101").Value = 695.0833
$ws.Range("K101").Value = 1696.5
$ws.Range("L101").Value = 2085.2499
$ws.Range("M101").Value = -74.5
$ws.Range("N101").Value = -5329.2499

# row 112
$ws.Range("H112").Value = 4498.75
$ws.Range("J112").Value = 4498.75
$ws.Range("L112").Value = 13496.25
$ws.Range("N112").Value = -15712.25

# row 121
$ws.Range("H121").Value = 491.33334
$ws.Range("J121").Value = 491.33334
$ws.Range("L121").Value = 1474.00002
$ws.Range("N121").Value = -4968.000019999999

# row 125
$ws.Range("H125").Value = 1330.7142
$ws.Range("I125").Value = 482
$ws.Range("K125").Value = 4338
$ws.Range("M125").Value = -1878

# row 137
$ws.Range("H137").Value = 2595.0417
$ws.Range("J137").Value = 2936.8125
$ws.Range("L137").Value = 8810.4375
$ws.Range("N137").Value = -13910.4375

# row 141
$ws.Range("H141").Value = 4033.5881
$ws.Range("I141").Value = 4583.857
$ws.Range("J141").Value = 1465.6666
$ws.Range("K141").Value = 13751.571
$ws.Range("L141").Value = 4396.9998
$ws.Range("M141").Value = -8571.571
$ws.Range("N141").Value = -14756.9998

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 5559.9883
$ws.Range("I32").Value = 4996.446
$ws.Range("J32").Value = 17253.5
$ws.Range("K32").Value = 4996.446
$ws.Range("L32").Value = 17253.5
$ws.Range("M32").Value = -4709.446
$ws.Range("N32").Value = -17827.5

# row 60
$ws.Range("H60").Value = 166743060
$ws.Range("J60").Value = 500050000
$ws.Range("L60").Value = 500050000
$ws.Range("N60").Value = -500051466

# row 61
$ws.Range("H61").Value = 4671810.5
$ws.Range("I61").Value = 5268499
$ws.Range("K61").Value = 5268499
$ws.Range("M61").Value = -5268287

# row 97
$ws.Range("H97").Value = 1358.875
$ws.Range("J97").Value = 2140.5557
$ws.Range("L97").Value = 2140.5557
$ws.Range("N97").Value = -3132.5557

# row 122
$ws.Range("H122").Value = 2101.276
$ws.Range("I122").Value = 1850.8846
$ws.Range("J122").Value = 4271.3335
$ws.Range("K122").Value = 5552.6538
$ws.Range("L122").Value = 12814.0005
$ws.Range("M122").Value = -3102.6538
$ws.Range("N122").Value = -17714.0005

# row 132
$ws.Range("H132").Value = 5266101.5
$ws.Range("I132").Value = 3544.7856
$ws.Range("J132").Value = 20001260
$ws.Range("K132").Value = 10634.3568
$ws.Range("L132").Value = 60003780
$ws.Range("M132").Value = -8104.356800000001
$ws.Range("N132").Value = -60008840

# row 136
$ws.Range("H136").Value = 4671810.5
$ws.Range("I136").Value = 5268499
$ws.Range("K136").Value = 15805497
$ws.Range("M136").Value = -15802947

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 8048
$ws.Range("I22").Value = 1274.2
$ws.Range("J22").Value = 16515.25
$ws.Range("K22").Value = 1274.2
$ws.Range("L22").Value = 16515.25
$ws.Range("M22").Value = -1101.2
$ws.Range("N22").Value = -16861.25

$ws = $wb.Worksheets.Item("CRP")
# row 105
$ws.Range("H105").Value = 5648.0625
$ws.Range("I105").Value = 1212.2
$ws.Range("K105").Value = 1212.2
$ws.Range("M105").Value = 534.8

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 3733.3333
$ws.Range("I80").Value = 3600
$ws.Range("K80").Value = 3600
$ws.Range("M80").Value = -2602

# row 83
$ws.Range("H83").Value = 3733.3333
$ws.Range("I83").Value = 3600
$ws.Range("K83").Value = 18000
$ws.Range("M83").Value = -13008

# row 97
$ws.Range("H97").Value = 2404.7646
$ws.Range("I97").Value = 2004.25
$ws.Range("J97").Value = 3366
$ws.Range("K97").Value = 2004.25
$ws.Range("L97").Value = 3366
$ws.Range("M97").Value = -1508.25
$ws.Range("N97").Value = -4358

$ws = $wb.Worksheets.Item("LTW")
# row 82
$ws.Range("H82").Value = 1885.4706
$ws.Range("I82").Value = 796.6
$ws.Range("K82").Value = 796.6
$ws.Range("M82").Value = -435.6

# row 85
$ws.Range("H85").Value = 1885.4706
$ws.Range("I85").Value = 796.6
$ws.Range("K85").Value = 796.6
$ws.Range("M85").Value = 451.4

# row 122
$ws.Range("H122").Value = 4455.3555
$ws.Range("J122").Value = 8562
$ws.Range("L122").Value = 25686
$ws.Range("N122").Value = -30586

# row 128
$ws.Range("H128").Value = 68332.336
$ws.Range("J128").Value = 68332.336
$ws.Range("L128").Value = 68332.336
$ws.Range("N128").Value = -78292.336

$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 1174.2858
$ws.Range("I81").Value = 1057.9166
$ws.Range("K81").Value = 2115.8332
$ws.Range("M81").Value = -1054.8332

# row 84
$ws.Range("H84").Value = 1174.2858
$ws.Range("I84").Value = 1057.9166
$ws.Range("K84").Value = 10579.166
$ws.Range("M84").Value = -5275.166000000001

# row 122
$ws.Range("H122").Value = 1752.9615
$ws.Range("I122").Value = 1399.125
$ws.Range("K122").Value = 4197.375
$ws.Range("M122").Value = -1747.375

# row 136
$ws.Range("H136").Value = 190437.03
$ws.Range("I136").Value = 1762.125
$ws.Range("K136").Value = 5286.375
$ws.Range("M136").Value = -2736.375

# row 141
$ws.Range("H141").Value = 94612.93
$ws.Range("J141").Value = 94612.93
$ws.Range("L141").Value = 94612.93
$ws.Range("N141").Value = -104972.93

